$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update B15: revised Title & Image Mismatch justification text, and drop its explicit style back to default ---
$ws.Range("B15").Value = 'The uploaded images do not accurately reflect the product description. Please ensure the images show the correct item count, packaging, and product details as stated in the reference name'
$ws.Range("B15").Style = "Normal"

# --- Append 4 new reason/justification rows (63-66) ---
$ws.Range("A63").Value = 'AI Generated Images'
$ws.Range("B63").Value = 'Please note that AI-generated images are not accepted, as they do not meet professional standards. We require images captured in a real, professional setting for approval.'
$ws.Range("A64").Value = 'Digitally Manipulated Images'
$ws.Range("B64").Value = 'The image has been rejected due to the use of digitally manipulated or composite elements. The duffle bag appears to be cut and pasted into the scene, resulting in unnatural shadows, inconsistent lighting, and a lack of realistic depth—indicative of photo editing rather than a genuine product shot.Please provide a real-life photograph that accurately reflects the product and its usage.'
$ws.Range("A65").Value = 'Nutrition Info Missing in Attribute'
$ws.Range("B65").Value = 'Since the nutritional information is currently displayed only on the image, it cannot be accepted. Please ensure all required details are fully filled out, properly mapped in the product attributes, and presented within the content/text fields—with correct structure and proper formatting. Information shown on images is not sufficient for compliance'
$ws.Range("A66").Value = 'Product Not Relased Yet'
$ws.Range("B66").Value = 'Product listing cannot be approved at this time as the item has not yet been officially launched or made available for public sale. we only accept submissions for products that are currently active, in-market, and available for purchase'

# Row 66's justification cell uses a new italic font style, matching the added cellXf/font in styles.xml
$ws.Range("B66").Font.Italic = $true

# --- Update the view: scroll near the bottom and select the next empty justification cell ---
$null = $ws.Range("B70").Select()
